$wb = $excel.ActiveWorkbook

# --- Sheet: Schedule ---
$ws1 = $wb.Worksheets.Item("Schedule")

$ws1.Range("A2").Value = 46055.16666666666
$ws1.Range("C2").Value = 12
$ws1.Range("D2").Value = 45.36
$ws1.Range("E2").Value = 474.7860487499999
$ws1.Range("F2").Value = 10.46706456679894
$ws1.Range("A3").Value = 46055.83333333334
$ws1.Range("C3").Value = 7.5
$ws1.Range("D3").Value = 28.35
$ws1.Range("E3").Value = 1083.7736715
$ws1.Range("F3").Value = 38.22834820105821
$ws1.Range("E4").Value = 648.43762425
$ws1.Range("F4").Value = 20.18168765172736

# --- Sheet: Detailed ---
$ws2 = $wb.Worksheets.Item("Detailed")

$ws2.Range("E6").Value = "OFF"
$ws2.Range("E7").Value = "OFF"
$ws2.Range("E8").Value = "OFF"
$ws2.Range("E9").Value = "OFF"
$ws2.Range("B38").Value = -6
$ws2.Range("B39").Value = 22.80121
$ws2.Range("B40").Value = 65.11238
$ws2.Range("C40").Value = "historical"
$ws2.Range("B41").Value = 71.40000000000001
$ws2.Range("C41").Value = "historical"
$ws2.Range("B42").Value = 73.19
$ws2.Range("C42").Value = "historical"
$ws2.Range("E42").Value = "ON"
$ws2.Range("B43").Value = 75.10932
$ws2.Range("C43").Value = "historical"
$ws2.Range("E43").Value = "ON"
$ws2.Range("B44").Value = 75.98278000000001
$ws2.Range("C44").Value = "historical"
$ws2.Range("E44").Value = "ON"
$ws2.Range("B45").Value = 73.19
$ws2.Range("C45").Value = "historical"
$ws2.Range("E45").Value = "ON"
$ws2.Range("B46").Value = 70.23766999999999
$ws2.Range("C46").Value = "historical"
$ws2.Range("B47").Value = 83.69302
$ws2.Range("C47").Value = "historical"
$ws2.Range("B48").Value = 84.79000000000001
$ws2.Range("C48").Value = "historical"
$ws2.Range("B49").Value = 73.17999
$ws2.Range("B50").Value = 66.31119
$ws2.Range("B51").Value = 79.60947
$ws2.Range("B52").Value = 74.92478
$ws2.Range("B53").Value = 73.20007
$ws2.Range("B54").Value = 67.21592
$ws2.Range("B55").Value = 67.72846
$ws2.Range("B56").Value = 73.20007
$ws2.Range("B58").Value = 77.93452000000001
$ws2.Range("B59").Value = 77.82303
$ws2.Range("B60").Value = 83.94917
$ws2.Range("B61").Value = 86.71747000000001
$ws2.Range("B62").Value = 87.42259
$ws2.Range("B63").Value = 97.5778
$ws2.Range("B64").Value = 84.79000000000001
$ws2.Range("B65").Value = 53.34919
$ws2.Range("B66").Value = 46.3042
$ws2.Range("B67").Value = 36.0601
$ws2.Range("B68").Value = 36.06005
$ws2.Range("B70").Value = 36.0601
$ws2.Range("B71").Value = 36.06049
$ws2.Range("B73").Value = 26.62484
$ws2.Range("B74").Value = 35.92515
$ws2.Range("B75").Value = 36.0601
$ws2.Range("B76").Value = 36.0601
$ws2.Range("B77").Value = 36.0601
$ws2.Range("B78").Value = 36.0601
$ws2.Range("B80").Value = 49.13934
$ws2.Range("B81").Value = 57.06007
$ws2.Range("B82").Value = 46.29534
$ws2.Range("B83").Value = 33.06049
$ws2.Range("B84").Value = 62.6256
$ws2.Range("B85").Value = 73.77512
$ws2.Range("B86").Value = 101.25
$ws2.Range("B87").Value = 100.91246
$ws2.Range("B88").Value = 105.79
$ws2.Range("B89").Value = 101.89163
$ws2.Range("B90").Value = 105.79
$ws2.Range("B91").Value = 106.53158
$ws2.Range("B92").Value = 109.41
$ws2.Range("B93").Value = 108.89
$ws2.Range("B94").Value = 101.33
$ws2.Range("B95").Value = 108.89
$ws2.Range("B96").Value = 120.62819
$ws2.Range("B97").Value = 108.89
